$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vars_meta_data")

# New variable metadata rows describing the Longhurst province fields and
# the Season field added to the dataset (3 new rows appended after row 40).
$rows = @(
    @{ Row = 41; ShortName = "Longhurst_Long";  LongName = "Longhurst province sample was taken in." },
    @{ Row = 42; ShortName = "Longhurst_Short"; LongName = "Longhurst province sample was taken in, shortened code." },
    @{ Row = 43; ShortName = "Season";          LongName = "Season sample was taken in." }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.ShortName   # var_short_name
    $ws.Cells.Item($row, 2).Value = $r.LongName    # var_long_name
    $ws.Cells.Item($row, 3).Value = "NA"           # var_sensor
    $ws.Cells.Item($row, 4).Value = "NA"           # var_unit
    $ws.Cells.Item($row, 5).Value = "Irregular"    # var_spatial_res
    $ws.Cells.Item($row, 6).Value = "Irregular"    # var_temporal_res
    $ws.Cells.Item($row, 7).Value = "Biology"      # var_discipline
    $ws.Cells.Item($row, 8).Value = 1              # visualize

    # Match the font styling (ArialMT 9pt) used by the neighboring rows'
    # A:D and G columns.
    $styled = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 4))
    $styled.Font.Name = "ArialMT"
    $styled.Font.Size = 9

    $discipline = $ws.Cells.Item($row, 7)
    $discipline.Font.Name = "ArialMT"
    $discipline.Font.Size = 9
}

# Restore selection/view state to match what was saved: range A41:H43
# selected with A41 active, no leftover scrolled topLeftCell.
$ws.Activate()
$ws.Range("A41:H43").Select()
